# Week 15 simulations - add new WR player row (E.Winston) to the WR sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# New player row with zeroed stat columns (B..J)
$ws.Cells.Item(11, 1).Value = "E.Winston"
for ($col = 2; $col -le 10; $col++) {
    $ws.Cells.Item(11, $col).Value = 0
}

# Make WR the active sheet and move the selection to J12, matching the
# author's saved view state after entering the new row.
$ws.Activate()
$ws.Range("J12").Select()
